# Update some hw and schedule
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# G2: Blog Intro/Bio entry wording tweak
$ws.Range("G2").Value = "Blog Intro/Bio (Due 8/28)`n* Blog RR plan  (Due 8/28)`n* Blog area of Add Health research interest (Due 8/28) `n* Personal Codebook/ Research Question Assignment  (Due 8/28)`n"

# G3: Personal codebook -> Blog about potential necessary recoding
$ws.Range("G3").Value = "Blog about potential necessary recoding  (Due Thu: 8/31) `n* Data management assignment (Due 9/4)`n* Citation Assignment (Due 9/7)"

# G4: Univariate graphing assignment due date change 9/14 -> 9/18
$ws.Range("G4").Value = "Univariate graphing assignment (Due 9/18)"

# Update the selected cell to G4
$ws.Range("G4").Select()
